$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 445
$ws.Cells.Item(2, 9).Value = 394.54544
$ws.Cells.Item(2, 10).Value = 1000
$ws.Cells.Item(2, 11).Value = 394.54544
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = -281.54544
$ws.Cells.Item(2, 14).Value = -1226
$ws.Cells.Item(11, 8).Value = 4868.0312
$ws.Cells.Item(11, 9).Value = 4868.0312
$ws.Cells.Item(11, 11).Value = 4868.0312
$ws.Cells.Item(11, 13).Value = -4728.0312
$ws.Cells.Item(39, 8).Value = 1168.2941
$ws.Cells.Item(39, 9).Value = 41.909092
$ws.Cells.Item(39, 10).Value = 3233.3333
$ws.Cells.Item(39, 11).Value = 125.727276
$ws.Cells.Item(39, 12).Value = 9699.999899999999
$ws.Cells.Item(39, 13).Value = 170.272724
$ws.Cells.Item(39, 14).Value = -10291.9999
$ws.Cells.Item(42, 8).Value = 126.2
$ws.Cells.Item(42, 9).Value = 126.2
$ws.Cells.Item(42, 11).Value = 378.6
$ws.Cells.Item(42, 13).Value = -148.6
$ws.Cells.Item(43, 8).Value = 3217.9092
$ws.Cells.Item(43, 9).Value = 3175
$ws.Cells.Item(43, 10).Value = 3242.4285
$ws.Cells.Item(43, 11).Value = 3175
$ws.Cells.Item(43, 12).Value = 3242.4285
$ws.Cells.Item(43, 13).Value = -3106
$ws.Cells.Item(43, 14).Value = -3380.4285
$ws.Cells.Item(94, 8).Value = 8318.556
$ws.Cells.Item(94, 9).Value = 9258.375
$ws.Cells.Item(94, 11).Value = 9258.375
$ws.Cells.Item(94, 13).Value = -8807.375
$ws.Cells.Item(103, 8).Value = 355.625
$ws.Cells.Item(103, 9).Value = 356.68182
$ws.Cells.Item(103, 10).Value = 344
$ws.Cells.Item(103, 11).Value = 1070.04546
$ws.Cells.Item(103, 12).Value = 1032
$ws.Cells.Item(103, 13).Value = -484.04546
$ws.Cells.Item(103, 14).Value = -2204
$ws.Cells.Item(132, 8).Value = 1532.8572
$ws.Cells.Item(132, 9).Value = 1465.8334
$ws.Cells.Item(132, 11).Value = 4397.5002
$ws.Cells.Item(132, 13).Value = -1867.5002

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 29373.195
$ws.Cells.Item(32, 9).Value = 29559.268
$ws.Cells.Item(32, 11).Value = 29559.268
$ws.Cells.Item(32, 13).Value = -29272.268
$ws.Cells.Item(122, 8).Value = 2422.923
$ws.Cells.Item(122, 9).Value = 3472.2856
$ws.Cells.Item(122, 10).Value = 1198.6666
$ws.Cells.Item(122, 11).Value = 10416.8568
$ws.Cells.Item(122, 12).Value = 3595.9998
$ws.Cells.Item(122, 13).Value = -7966.856800000001
$ws.Cells.Item(122, 14).Value = -8495.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 226.33333
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 12).Value = 150
$ws.Cells.Item(4, 14).Value = -380
$ws.Cells.Item(86, 8).Value = 106786.9
$ws.Cells.Item(86, 9).Value = 1414.4166
$ws.Cells.Item(86, 11).Value = 1414.4166
$ws.Cells.Item(86, 13).Value = -291.4166
$ws.Cells.Item(89, 8).Value = 106786.9
$ws.Cells.Item(89, 9).Value = 1414.4166
$ws.Cells.Item(89, 11).Value = 7072.083000000001
$ws.Cells.Item(89, 13).Value = -1456.083000000001
$ws.Cells.Item(94, 8).Value = 2032.7368
$ws.Cells.Item(94, 9).Value = 2227.0833
$ws.Cells.Item(94, 11).Value = 2227.0833
$ws.Cells.Item(94, 13).Value = -1776.0833
$ws.Cells.Item(105, 8).Value = 31260062
$ws.Cells.Item(105, 9).Value = 50014100
$ws.Cells.Item(105, 10).Value = 3332.25
$ws.Cells.Item(105, 11).Value = 50014100
$ws.Cells.Item(105, 12).Value = 3332.25
$ws.Cells.Item(105, 13).Value = -50012353
$ws.Cells.Item(105, 14).Value = -6826.25
$ws.Cells.Item(119, 8).Value = 0
$ws.Cells.Item(119, 10).Value = 0
$ws.Cells.Item(119, 12).Value = 0
$ws.Cells.Item(125, 8).Value = 34259.668
$ws.Cells.Item(125, 10).Value = 34259.668
$ws.Cells.Item(125, 12).Value = 34259.668
$ws.Cells.Item(125, 14).Value = -44099.668
$ws.Cells.Item(126, 8).Value = 0
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(119, 14).ClearContents()
$ws.Cells.Item(126, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 238.53334
$ws.Cells.Item(7, 9).Value = 262.9
$ws.Cells.Item(7, 10).Value = 189.8
$ws.Cells.Item(7, 11).Value = 262.9
$ws.Cells.Item(7, 12).Value = 189.8
$ws.Cells.Item(7, 13).Value = -149.9
$ws.Cells.Item(7, 14).Value = -415.8
$ws.Cells.Item(16, 8).Value = 2052.25
$ws.Cells.Item(16, 9).Value = 1936.6666
$ws.Cells.Item(16, 11).Value = 1936.6666
$ws.Cells.Item(16, 13).Value = -1649.6666
$ws.Cells.Item(31, 8).Value = 2873.3489
$ws.Cells.Item(31, 9).Value = 1699.4231
$ws.Cells.Item(31, 11).Value = 1699.4231
$ws.Cells.Item(31, 13).Value = -1404.4231
$ws.Cells.Item(34, 8).Value = 2873.3489
$ws.Cells.Item(34, 9).Value = 1699.4231
$ws.Cells.Item(34, 11).Value = 1699.4231
$ws.Cells.Item(34, 13).Value = -1497.4231
$ws.Cells.Item(68, 8).Value = 599999.7
$ws.Cells.Item(68, 10).Value = 999999
$ws.Cells.Item(68, 12).Value = 999999
$ws.Cells.Item(68, 14).Value = -1001497
$ws.Cells.Item(71, 8).Value = 599999.7
$ws.Cells.Item(71, 10).Value = 999999
$ws.Cells.Item(71, 12).Value = 2999997
$ws.Cells.Item(71, 14).Value = -3007485
$ws.Cells.Item(113, 8).Value = 2052.25
$ws.Cells.Item(113, 9).Value = 1936.6666
$ws.Cells.Item(113, 11).Value = 1936.6666
$ws.Cells.Item(113, 13).Value = 233.3334
$ws.Cells.Item(132, 8).Value = 20911.973
$ws.Cells.Item(132, 10).Value = 54782.918
$ws.Cells.Item(132, 12).Value = 164348.754
$ws.Cells.Item(132, 14).Value = -169408.754
$ws.Cells.Item(134, 8).Value = 4540.3516
$ws.Cells.Item(134, 9).Value = 4111.387
$ws.Cells.Item(134, 11).Value = 12334.161
$ws.Cells.Item(134, 13).Value = -9799.161

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 502.14285
$ws.Cells.Item(22, 9).Value = 215
$ws.Cells.Item(22, 10).Value = 1019
$ws.Cells.Item(22, 11).Value = 645
$ws.Cells.Item(22, 12).Value = 3057
$ws.Cells.Item(22, 13).Value = -476
$ws.Cells.Item(22, 14).Value = -3395
$ws.Cells.Item(27, 8).Value = 502.14285
$ws.Cells.Item(27, 9).Value = 215
$ws.Cells.Item(27, 10).Value = 1019
$ws.Cells.Item(27, 11).Value = 645
$ws.Cells.Item(27, 12).Value = 3057
$ws.Cells.Item(27, 13).Value = -543
$ws.Cells.Item(27, 14).Value = -3261
$ws.Cells.Item(28, 8).Value = 2443.7144
$ws.Cells.Item(28, 10).Value = 3894
$ws.Cells.Item(28, 12).Value = 11682
$ws.Cells.Item(28, 14).Value = -12146
$ws.Cells.Item(58, 8).Value = 2766.258
$ws.Cells.Item(58, 9).Value = 1627.5
$ws.Cells.Item(58, 10).Value = 2844.7932
$ws.Cells.Item(58, 11).Value = 4882.5
$ws.Cells.Item(58, 12).Value = 8534.3796
$ws.Cells.Item(58, 13).Value = -4754.5
$ws.Cells.Item(58, 14).Value = -8790.3796
$ws.Cells.Item(75, 8).Value = 1553.75
$ws.Cells.Item(75, 10).Value = 2395
$ws.Cells.Item(75, 12).Value = 7185
$ws.Cells.Item(75, 14).Value = -9181
$ws.Cells.Item(78, 8).Value = 1553.75
$ws.Cells.Item(78, 10).Value = 2395
$ws.Cells.Item(78, 12).Value = 21555
$ws.Cells.Item(78, 14).Value = -31539
$ws.Cells.Item(113, 8).Value = 2020.4
$ws.Cells.Item(113, 9).Value = 1535.6666
$ws.Cells.Item(113, 10).Value = 2105.9412
$ws.Cells.Item(113, 11).Value = 4606.9998
$ws.Cells.Item(113, 12).Value = 6317.823600000001
$ws.Cells.Item(113, 13).Value = -2436.9998
$ws.Cells.Item(113, 14).Value = -10657.8236
$ws.Cells.Item(131, 8).Value = 15155540
$ws.Cells.Item(131, 9).Value = 37037944
$ws.Cells.Item(131, 10).Value = 6184.3076
$ws.Cells.Item(131, 11).Value = 111113832
$ws.Cells.Item(131, 12).Value = 18552.9228
$ws.Cells.Item(131, 13).Value = -111108792
$ws.Cells.Item(131, 14).Value = -28632.9228
$ws.Cells.Item(134, 8).Value = 2594
$ws.Cells.Item(134, 9).Value = 2594
$ws.Cells.Item(134, 11).Value = 7782
$ws.Cells.Item(134, 13).Value = -2712

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(32, 8).Value = 119250
$ws.Cells.Item(32, 10).Value = 120000
$ws.Cells.Item(32, 12).Value = 120000
$ws.Cells.Item(32, 14).Value = -120592
$ws.Cells.Item(113, 8).Value = 2964.9412
$ws.Cells.Item(113, 9).Value = 2892
$ws.Cells.Item(113, 10).Value = 3098.6667
$ws.Cells.Item(113, 11).Value = 2892
$ws.Cells.Item(113, 12).Value = 3098.6667
$ws.Cells.Item(113, 13).Value = -722
$ws.Cells.Item(113, 14).Value = -7438.6667
$ws.Cells.Item(138, 8).Value = 80000
$ws.Cells.Item(138, 10).Value = 80000
$ws.Cells.Item(138, 12).Value = 80000
$ws.Cells.Item(138, 14).Value = -90280

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 1520
$ws.Cells.Item(13, 9).Value = 790
$ws.Cells.Item(13, 10).Value = 5900
$ws.Cells.Item(13, 11).Value = 790
$ws.Cells.Item(13, 12).Value = 5900
$ws.Cells.Item(13, 13).Value = -650
$ws.Cells.Item(13, 14).Value = -6180
$ws.Cells.Item(25, 8).Value = 8266.611000000001
$ws.Cells.Item(25, 10).Value = 9181.25
$ws.Cells.Item(25, 12).Value = 9181.25
$ws.Cells.Item(25, 14).Value = -9641.25
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 10).Value = 0
$ws.Cells.Item(33, 12).Value = 0
$ws.Cells.Item(124, 8).Value = 76951.39999999999
$ws.Cells.Item(124, 10).Value = 76951.39999999999
$ws.Cells.Item(124, 12).Value = 76951.39999999999
$ws.Cells.Item(124, 14).Value = -86771.39999999999
$ws.Cells.Item(33, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 8552.6
$ws.Cells.Item(132, 9).Value = 7688.25
$ws.Cells.Item(132, 11).Value = 23064.75
$ws.Cells.Item(132, 13).Value = -20534.75
$ws.Cells.Item(136, 8).Value = 3324352.2
$ws.Cells.Item(136, 9).Value = 4329926
$ws.Cells.Item(136, 10).Value = 5958.5
$ws.Cells.Item(136, 11).Value = 12989778
$ws.Cells.Item(136, 12).Value = 17875.5
$ws.Cells.Item(136, 13).Value = -12987228
$ws.Cells.Item(136, 14).Value = -22975.5
